$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "D2" "27.357.39"
Set-TextValue "E2" "  -1.98%  "
Set-TextValue "D3" "1.652.55"
Set-TextValue "E3" "  -0.72%  "
Set-TextValue "E4" "  -0.33%  "
Set-TextValue "D5" "213.34"
Set-TextValue "E5" "  -0.96%  "
Set-TextValue "E6" "  -0.70%  "
Set-TextValue "E7" "  -0.29%  "
Set-TextValue "D8" "23.93"
Set-TextValue "E8" "  +1.90%  "
Set-TextValue "D9" "0.260"
Set-TextValue "E9" "  -0.66%  "
Set-TextValue "E10" "  -0.85%  "
Set-TextValue "D11" "0.0877"
Set-TextValue "E11" "  -0.51%  "
Set-TextValue "D12" "1.887.41"
Set-TextValue "E12" "  -0.61%  "
Set-TextValue "D13" "1.650.41"
Set-TextValue "E13" "  -1.08%  "
Set-TextValue "E14" "  -1.82%  "
Set-TextValue "D15" "0.569"
Set-TextValue "E15" "  +3.38%  "
Set-TextValue "D16" "65.60"
Set-TextValue "E16" "  -1.14%  "
Set-TextValue "D17" "27.374.04"
Set-TextValue "E17" "  -1.76%  "
Set-TextValue "D18" "231.99"
Set-TextValue "E18" "  -6.60%  "
Set-TextValue "D19" "0.0₃0725"
Set-TextValue "E19" "  -1.31%  "
Set-TextValue "D20" "7.50"
Set-TextValue "E20" "  -1.13%  "
Set-TextValue "E21" "  -0.50%  "
Set-TextValue "D22" "4.36"
Set-TextValue "E22" "  -2.52%  "
Set-TextValue "E23" "  -1.33%  "
Set-TextValue "E24" "  -1.10%  "
Set-TextValue "D25" "146.97"
Set-TextValue "E25" "  +0.05%  "
Set-TextValue "D26" "7.16"
Set-TextValue "E26" "  -1.16%  "
Set-TextValue "D27" "15.86"
Set-TextValue "E27" "  -2.37%  "
Set-TextValue "E28" "  -0.39%  "
Set-TextValue "E29" "  -0.89%  "
Set-TextValue "D30" "0.0498"
Set-TextValue "E30" "  -0.40%  "
Set-TextValue "E31" "  -3.84%  "
Set-TextValue "E32" "  -1.33%  "
Set-TextValue "D33" "1.459.98"
Set-TextValue "E33" "  +3.29%  "
Set-TextValue "E34" "  -0.95%  "
Set-TextValue "D35" "1.54"
Set-TextValue "E35" "  -1.71%  "
Set-TextValue "E36" "  -0.30%  "
Set-TextValue "D37" "0.907"
Set-TextValue "E37" "  -2.37%  "
Set-TextValue "D38" "0.570"
Set-TextValue "E38" "  -1.65%  "
Set-TextValue "E39" "  -0.39%  "
Set-TextValue "E40" "  +1.09%  "
Set-TextValue "E41" "  -0.40%  "
Set-TextValue "E42" "  +0.10%  "
Set-TextValue "D43" "65.20"
Set-TextValue "E43" "  -5.90%  "
Set-TextValue "E44" "  -0.71%  "
Set-TextValue "D45" "1.795.48"
Set-TextValue "E45" "  -0.81%  "
Set-TextValue "E46" "  -0.98%  "
Set-TextValue "D48" "88.25"
Set-TextValue "E48" "  -0.15%  "
Set-TextValue "D49" "0.0₆0106"
Set-TextValue "E49" "  -1.23%  "
Set-TextValue "E50" "  -0.47%  "
Set-TextValue "D51" "7.73"
Set-TextValue "E51" "  -0.27%  "
